# 4 & 5 Modification 2
# The "Уведомление о доставках" (courier notification) use-case row is removed
# from the use-case table. That row was row 12 on the sheet; deleting it shifts
# all subsequent rows up by one and drops the now-unused shared strings
# belonging to that row (handled automatically by the engine on save).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Select the entire row first so the resulting selection/active cell lands on
# row 12 (the row that used to be row 13), matching how Excel behaves when a
# user selects a row header and deletes it.
$ws.Rows(12).EntireRow.Select()
$ws.Rows(12).Delete()

# Best-effort: scroll the view so row 13 becomes the first visible row.
$excel.ActiveWindow.ScrollRow = 13
$excel.ActiveWindow.ScrollColumn = 1
